$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("car inventory")

# Set column N width (like the bestFit/customWidth columns already present)
$ws.Columns.Item(14).ColumnWidth = 16.5

# Row 2 formulas (anchors for the shared formula ranges below)
$ws.Range("M2").Formula = '=IF(H2<=L2,"Y", "N")'
$ws.Range("N2").Formula = '=CONCATENATE(B2,F2,D2,UPPER(LEFT(J2,3)),RIGHT(A2,3))'

# Rows 3 through 53: same formulas, each row referencing its own row number
for ($r = 3; $r -le 53; $r++) {
    $ws.Range("M$r").Formula = "=IF(H$r<=L$r,`"Y`", `"N`")"
    $ws.Range("N$r").Formula = "=CONCATENATE(B$r,F$r,D$r,UPPER(LEFT(J$r,3)),RIGHT(A$r,3))"
}

# Update the selected cell to match the final state
$ws.Range("N53").Select()
